$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume refresh (scheduled GitHub Actions update).
# Force Text number format on D (Price) / E (Volume 1h %) cells before
# assigning, so values like "15.80" or "0.0000110" keep their exact
# textual representation instead of being normalized as numbers.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.632.97'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.14%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.269.33'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.20%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '119.49'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +6.34%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '269.48'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +1.78%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +3.52%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.36%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.621'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +2.28%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '47.44'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.61%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0944'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.43%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +6.43%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.36%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.80'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +2.27%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.915'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +7.29%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.612.22'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.17%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.266.02'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.34%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.537.29'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.87%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0000110'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.65%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.91'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +2.39%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.46'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.70%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.39'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -4.26%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.03'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +4.74%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '234.59'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.30%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.60'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.68%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.27'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +8.81%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +1.62%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '41.71'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +3.69%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.94%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.38%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '174.84'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +1.63%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '21.53'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.12%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0917'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +1.49%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.72'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.68%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.131'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +3.08%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +12.14%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0379'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +7.96%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.67'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.62%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +4.80%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.21%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '13.67'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.40%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +1.68%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '71.99'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -6.32%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.999'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.23%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.32%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.71'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -6.26%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '75.40'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +36.87%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.28'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +2.53%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +18.87%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +2.02%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.59'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.17%  '
